# Applies the Jul 27 2023 cryptos-list update (GitHub Actions scrape).
# Updates Price (D) / Volume(1h) (E) columns for each coin row, and
# fixes the row order for WrappedEther / Polygon / Polkadot (rows 12-14).
# Numeric-looking price strings are apostrophe-prefixed so Excel keeps
# them as literal text (preserving trailing zeros / exact formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.374.46"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.873.73"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'0.7133"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").Value = "'239.56"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.07810"
$ws.Range("E8").Value = "  -2.24%  "
$ws.Range("D9").Value = "'0.3078"
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("D10").Value = "'24.82"
$ws.Range("E10").Value = "  +5.98%  "
$ws.Range("D11").Value = "'0.08246"
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.867.74"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.7236"
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.265"
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("D15").Value = "'91.37"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "29.442.74"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").Value = "'5.883"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").Value = "'0.000007916"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").Value = "'243.69"
$ws.Range("E19").Value = "  +3.39%  "
$ws.Range("D20").Value = "'13.27"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "'0.9992"
$ws.Range("D22").Value = "'7.945"
$ws.Range("E22").Value = "  +7.74%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'0.1548"
$ws.Range("E24").Value = "  +8.15%  "
$ws.Range("D25").Value = "'163.42"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("D26").Value = "'8.980"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").Value = "'18.29"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("E28").Value = "  -4.29%  "
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").Value = "'4.364"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "'4.109"
$ws.Range("E31").Value = "  +2.52%  "
$ws.Range("D32").Value = "'0.05266"
$ws.Range("D33").Value = "'1.920"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").Value = "'1.197"
$ws.Range("E34").Value = "  +3.40%  "
$ws.Range("D35").Value = "'0.7177"
$ws.Range("E35").Value = "  +2.91%  "
$ws.Range("D36").Value = "'2.679"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").Value = "'0.01858"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").Value = "1.206.49"
$ws.Range("E38").Value = "  +7.50%  "
$ws.Range("D39").Value = "'2.708"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "'0.9061"
$ws.Range("E40").Value = "  -2.77%  "
$ws.Range("D41").Value = "'6.054"
$ws.Range("E41").Value = "  +3.94%  "
$ws.Range("E42").Value = "  +4.07%  "
$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "'103.21"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "'0.5340"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("E46").Value = "  +4.30%  "
$ws.Range("D47").Value = "'1.755"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "'2.897"
$ws.Range("E48").Value = "  +10.15%  "
$ws.Range("D49").Value = "'0.4307"
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("D50").Value = "'9.225"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("D51").Value = "'0.9995"
$ws.Range("E51").Value = "  +0.08%  "
